# "create migration for setup": relabel the driver-format header row with
# friendly display names, drop the old "phone" header (column G keeps its
# place/format but is left blank), and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Address"
$ws.Range("C1").Value = "Salary"
$ws.Range("D1").Value = "License No"
$ws.Range("E1").Value = "License Eexp"
$ws.Range("F1").Value = "Phone"

# G1 ("phone") loses its header text but the cell itself (and its format)
# stays in the used range, same as Excel leaves behind after clearing text
# while keeping formatting.
$ws.Range("G1").ClearContents()
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J8").Select()
